$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "T"
$ws.Range("D3").Value = "T"
$ws.Range("D4").Value = "T"
$ws.Range("D5").Value = "T"
$ws.Range("D6").Value = "T"
$ws.Range("D7").Value = "T"
$ws.Range("D8").Value = "T"
$ws.Range("D9").Value = "T"
$ws.Range("D10").Value = "T"
$ws.Range("D11").Value = "T"
$ws.Range("D12").Value = "T"
$ws.Range("D13").Value = "T"
$ws.Range("D14").Value = "T"
$ws.Range("D15").Value = "T"
$ws.Range("D16").Value = "T"
$ws.Range("D17").Value = "T"
$ws.Range("D18").Value = "T"
$ws.Range("D19").Value = "T"
$ws.Range("D20").Value = "T"
$ws.Range("D21").Value = "T"
$ws.Range("D22").Value = "T"
$ws.Range("D23").Value = "T"
$ws.Range("D24").Value = "S"
$ws.Range("D25").Value = "S"
$ws.Range("D26").Value = "S"
$ws.Range("D27").Value = "S"
$ws.Range("D28").Value = "SN"
$ws.Range("D29").Value = "S"
$ws.Range("D30").Value = "SN"
$ws.Range("D31").Value = "S"
$ws.Range("D32").Value = "SN"
$ws.Range("D33").Value = "SN"
$ws.Range("D34").Value = "S"
$ws.Range("D35").Value = "SN"
$ws.Range("G35").Value = "2 - Relating to Another S"
$ws.Range("D36").Value = "T"
$ws.Range("D37").Value = "S"
$ws.Range("D38").Value = "SN"
$ws.Range("G38").Value = "2 - Relating to Another S"
$ws.Range("D39").Value = "SN"
$ws.Range("G39").Value = "2 - Relating to Another S"
$ws.Range("D40").Value = "SN"
$ws.Range("G40").Value = "2 - Relating to Another S"
$ws.Range("D41").Value = "S"
$ws.Range("D42").Value = "SN"
$ws.Range("D43").Value = "SN"
$ws.Range("G43").Value = "2 - Relating to Another S"
$ws.Range("D44").Value = "SN"
$ws.Range("D45").Value = "SN"
$ws.Range("G45").Value = "2 - Relating to Another S"
$ws.Range("D46").Value = "SN"
$ws.Range("D47").Value = "SN"
$ws.Range("G47").Value = "2 - Relating to Another S"
$ws.Range("D48").Value = "SN"
$ws.Range("G48").Value = "2 - Relating to Another S"
$ws.Range("D49").Value = "S"
$ws.Range("D50").Value = "SN"
$ws.Range("D51").Value = "SN"
$ws.Range("G51").Value = "2 - Relating to Another S"
$ws.Range("D52").Value = "S"
$ws.Range("D53").Value = "S"
$ws.Range("D54").Value = "S"
$ws.Range("D55").Value = "S"
$ws.Range("D56").Value = "SN"
$ws.Range("D57").Value = "T"
$ws.Range("D58").Value = "T"
$ws.Range("D59").Value = "T"
$ws.Range("D60").Value = "T"
$ws.Range("D61").Value = "T"
$ws.Range("D62").Value = "S"
$ws.Range("D63").Value = "SN"
$ws.Range("G63").Value = "2 - Relating to Another S"
$ws.Range("D64").Value = "S"
$ws.Range("D65").Value = "SN"
$ws.Range("G65").Value = "2 - Relating to Another S"
$ws.Range("D66").Value = "S"
$ws.Range("D67").Value = "SN"
$ws.Range("D68").Value = "S"
$ws.Range("G68").Value = "3 - Asking for Information"
$ws.Range("D69").Value = "S"
$ws.Range("D70").Value = "S"
$ws.Range("D71").Value = "SN"
$ws.Range("G71").Value = "2 - Relating to Another S"
$ws.Range("D72").Value = "SN"
$ws.Range("G72").Value = "2 - Relating to Another S"
$ws.Range("D73").Value = "S"
$ws.Range("D74").Value = "S"
$ws.Range("D75").Value = "SN"
$ws.Range("G75").Value = "2 - Relating to Another S"
$ws.Range("D76").Value = "S"
$ws.Range("D77").Value = "S"
$ws.Range("D78").Value = "S"
$ws.Range("D79").Value = "S"
$ws.Range("D80").Value = "S"
$ws.Range("D81").Value = "SN"
$ws.Range("G81").Value = "2 - Relating to Another S"
$ws.Range("D82").Value = "S"
$ws.Range("D83").Value = "S"
$ws.Range("D84").Value = "SN"
$ws.Range("G84").Value = "2 - Relating to Another S"
$ws.Range("D85").Value = "S"
$ws.Range("D86").Value = "S"
$ws.Range("D87").Value = "SN"
$ws.Range("G87").Value = "3 - Asking for Information"
$ws.Range("D88").Value = "SN"
$ws.Range("G88").Value = "2 - Relating to Another S"
$ws.Range("D89").Value = "SN"
$ws.Range("D90").Value = "S"
$ws.Range("D91").Value = "SN"
$ws.Range("D92").Value = "SN"
$ws.Range("D93").Value = "T"
$ws.Range("D94").Value = "T"
$ws.Range("D95").Value = "T"
$ws.Range("D96").Value = "T"
$ws.Range("D97").Value = "S"
$ws.Range("D98").Value = "T"
$ws.Range("D99").Value = "T"
$ws.Range("D100").Value = "T"
$ws.Range("D101").Value = "T"
$ws.Range("D102").Value = "T"
$ws.Range("D103").Value = "T"
$ws.Range("D104").Value = "T"
$ws.Range("D105").Value = "T"
$ws.Range("D106").Value = "T"
$ws.Range("D107").Value = "T"
$ws.Range("D108").Value = "S"
$ws.Range("D109").Value = "S"
$ws.Range("D110").Value = "S"
$ws.Range("D111").Value = "S"
$ws.Range("D112").Value = "S"
$ws.Range("D113").Value = "S"
$ws.Range("D114").Value = "S"
$ws.Range("D115").Value = "S"
$ws.Range("D116").Value = "S"
$ws.Range("D117").Value = "S"
$ws.Range("D118").Value = "S"
$ws.Range("G118").Value = "2 - Relating to Another S"
$ws.Range("D119").Value = "S"
$ws.Range("D120").Value = "S"
$ws.Range("D121").Value = "S"
$ws.Range("D122").Value = "S"
